$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = 2158.5577258019043
$ws.Range("D7").Value = 2224.8751303977206

$ws.Range("B11").Value = 842.60026661799202

$ws.Range("B14").Value = 975.42889992313644
$ws.Range("C14").Value = 52513.94549326141

$ws.Range("B17").Value = 3210.5350387796743
$ws.Range("C17").Value = 109817.62953418112
$ws.Range("D17").Value = 2582.0413660436884

$ws.Range("B18").Value = 3160.9167526535193
$ws.Range("C18").Value = 154852.87741765584
$ws.Range("D18").Value = 5585.0394878110656

$ws.Range("B19").Value = 812.19374298080118
$ws.Range("C19").Value = 231550.00333940168

$ws.Range("B20").Value = 6888.0287740829299
$ws.Range("C20").Value = 270181.46518217353
$ws.Range("D20").Value = 7466.419295420651

$ws.Range("B22").Value = 1053.0900335721303

$ws.Range("B23").Value = 581.02458859681474
$ws.Range("C23").Value = 77712.786054656841

$ws.Range("B24").Value = 573.39782350605412

$ws.Range("B25").Value = 600.64919257420979

$ws.Range("B27").Value = 3327.9327645134244
$ws.Range("C27").Value = 239064.21059067981

$ws.Range("B28").Value = 7419.3381679044796
$ws.Range("C28").Value = 287053.8356120314
$ws.Range("D28").Value = 6516.8801523460479
